$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.189719438552856
$ws.Range("B1").Value = 1.806123971939087
$ws.Range("C1").Value = 6.611209392547607
$ws.Range("D1").Value = 2.284949064254761
$ws.Range("E1").Value = 1.192989826202393
